# Payment order template: add the {CURR_DATE} field right after the city
# name "Τρίπολη" and before the trailing comma, e.g.
#   "Τρίπολη,"  ->  "Τρίπολη {CURR_DATE},"

$d = $word.ActiveDocument

# Locate the run of text "Τρίπολη," (city name immediately followed by a
# comma) - this is the signature/date line near the bottom of the document.
$rng = $d.Content
$found = $rng.Find.Execute("Τρίπολη,", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target text 'Τρίπολη,'"
}

# Strip the trailing comma from the found text, leaving just "Τρίπολη".
$rng.Text = "Τρίπολη"

# Move to right after "Τρίπολη" and type the new content in one shot:
# a space, the literal "{CURR_DATE}" placeholder (to be bolded), and the
# comma that used to directly follow the city name.
$rng.Collapse(0)
$rng.InsertAfter(" {CURR_DATE},")
$s = $rng.Start

# Re-apply formatting per character range so each logical piece keeps its
# own run in the saved XML (mirrors how Word keeps separately-typed /
# separately-formatted spans as distinct <w:r> elements instead of always
# silently coalescing runs that happen to end up with equal formatting).

# " " (space) - plain; flip Bold on/off to force a run boundary against the
# preceding plain "Τρίπολη" text.
$spaceRng = $d.Range($s, $s + 1)
$spaceRng.Bold = 1
$spaceRng.Bold = 0

# "{" - bold
$d.Range($s + 1, $s + 2).Bold = 1

# "CURR_DATE" - bold
$d.Range($s + 2, $s + 11).Bold = 1

# "}" - bold
$d.Range($s + 11, $s + 12).Bold = 1

# "," - plain; flip Bold on/off so it doesn't silently re-merge with the
# bold "}" before it nor swallow the untouched text that follows it.
$commaRng = $d.Range($s + 12, $s + 13)
$commaRng.Bold = 1
$commaRng.Bold = 0

Write-Output "Inserted {CURR_DATE} after Τρίπολη"
